$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 265, shifting existing rows 265:289 down to 266:290.
$ws.Rows.Item(265).Insert()

# Populate the newly-inserted row 265 with the new weekly Granada (Wonderfull / Primera) record.
$ws.Range("A265").Value = 10
$ws.Range("B265").Value = "Vega Modelo de Temuco"
$ws.Range("C265").Value = "La Araucanía"
$ws.Range("D265").Value = 45166
$ws.Range("E265").Value = 9
$ws.Range("F265").Value = "Fruta"
$ws.Range("G265").Value = 100104
$ws.Range("H265").Value = "Frutos de pepita"
$ws.Range("I265").Value = 100104001
$ws.Range("J265").Value = "Granada"
$ws.Range("K265").Value = "Wonderfull"
$ws.Range("L265").Value = "Primera"
$ws.Range("M265").Value = 55
$ws.Range("N265").Value = 15000
$ws.Range("O265").Value = 15000
$ws.Range("P265").Value = 15000
$ws.Range("Q265").Value = "$/bandeja 10 kilos granel"
$ws.Range("R265").Value = "Provincia de Limarí"
$ws.Range("S265").Value = 1500
$ws.Range("T265").Value = 10
